$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data for new column L (2023)
$ws.Range("L4").Value = 2023
$ws.Range("L5").Value = 21.9
$ws.Range("L7").Value = 23.2
$ws.Range("L8").Value = 19.9
$ws.Range("L10").Value = 17.1
$ws.Range("L11").Value = 14.7
$ws.Range("L12").Value = 16.9
$ws.Range("L13").Value = 29.1
$ws.Range("L14").Value = 28.3
$ws.Range("L15").Value = 33.4
$ws.Range("L16").Value = 22.1
$ws.Range("L17").Value = 16.6
$ws.Range("L18").Value = 14.7
$ws.Range("L19").Value = 9.1
$ws.Range("L20").Value = 12.5
$ws.Range("L21").Value = 11.6
$ws.Range("L23").Value = 22.6
$ws.Range("L24").Value = 21.5
$ws.Range("L26").Value = 24.3
$ws.Range("L27").Value = 44.96
$ws.Range("L28").Value = 22.5
$ws.Range("L29").Value = 17.4
$ws.Range("L30").Value = 30
$ws.Range("L31").Value = 25.3
$ws.Range("L32").Value = 21
$ws.Range("L33").Value = 22
$ws.Range("L34").Value = 15
$ws.Range("L35").Value = 6.3
$ws.Range("L36").Value = 9
$ws.Range("L37").Value = 39.2
$ws.Range("L38").Value = 12
$ws.Range("L39").Value = 41.3
$ws.Range("L40").Value = 17.4
$ws.Range("L41").Value = 25.6
$ws.Range("L42").Value = 18.9
$ws.Range("L43").Value = 15.3
$ws.Range("L44").Value = 18.1
$ws.Range("L45").Value = 53
$ws.Range("L46").Value = 5.2

# Empty (but styled) cells in column L for "header" rows
$emptyRows = @(6, 9, 22, 25)
foreach ($r in $emptyRows) {
    $ws.Range("L$r").Value = $null
}

# Copy styles from column K to column L for rows 4-46
for ($r = 4; $r -le 46; $r++) {
    $ws.Range("K$r").Copy() | Out-Null
    $ws.Range("L$r").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}
$excel.CutCopyMode = 0

# Remove the selection that was set to L12 in the original file
$ws.Range("A1").Select() | Out-Null

$wb.Save()
